$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 3 new rows of test data (rows 9, 10, 11) - "1 extra test" entries
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 1000
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 0.2
$ws.Range("E9").Value = 100
$ws.Range("F9").Value = 100
$ws.Range("G9").Value = 99

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 100
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 0.1
$ws.Range("E10").Value = 100
$ws.Range("F10").Value = 100
$ws.Range("G10").Value = 59

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 100
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 0.2
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 100
$ws.Range("G11").Value = 56

# Update the defined name range to include the new rows
$wb.Names.Item("neuralNetworkLog").RefersTo = "=Sheet1!`$A`$1:`$G`$11"

# Update the chart's source data series formula to include new rows
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(Sheet1!`$G`$1,,Sheet1!`$G`$2:`$G`$11,1)"

# Reposition/resize the chart to its new anchor location
$co.Left = 554.810546875
$co.Top = 16.12496062992126
$co.Width = 433.0625
$co.Height = 216

# Update the active selection
$ws.Range("R5").Select()

$wb.Save()
